$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Volume/percentage cells (E column) are always text with padding spaces,
# and Price cells (D column) must stay as text too (e.g. "210.60", "1.00").
# Force text format so Excel does not auto-convert numeric-looking strings.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.056.13"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").Value = "1.563.30"
$ws.Range("E3").Value = "  +0.58%  "

$ws.Range("E4").Value = "  +0.48%  "

$ws.Range("D5").Value = "210.60"
$ws.Range("E5").Value = "  +1.79%  "

$ws.Range("D6").Value = "0.490"
$ws.Range("E6").Value = "  +0.32%  "

$ws.Range("E7").Value = "  +0.51%  "

$ws.Range("D8").Value = "21.90"
$ws.Range("E8").Value = "  -0.58%  "

$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").Value = "0.0597"
$ws.Range("E10").Value = "  +0.34%  "

$ws.Range("D11").Value = "0.0862"
$ws.Range("E11").Value = "  +0.63%  "

$ws.Range("D12").Value = "1.782.94"
$ws.Range("E12").Value = "  +0.50%  "

$ws.Range("D13").Value = "1.559.86"
$ws.Range("E13").Value = "  +0.81%  "

$ws.Range("D14").Value = "3.77"
$ws.Range("E14").Value = "  +0.29%  "

$ws.Range("D15").Value = "0.518"
$ws.Range("E15").Value = "  -0.34%  "

$ws.Range("D16").Value = "27.062.37"
$ws.Range("E16").Value = "  +0.52%  "

$ws.Range("D17").Value = "61.97"
$ws.Range("E17").Value = "  +0.51%  "

$ws.Range("E18").Value = "  -0.81%  "

$ws.Range("D19").Value = "215.11"
$ws.Range("E19").Value = "  -1.04%  "

$ws.Range("D20").Value = "7.36"
$ws.Range("E20").Value = "  +0.72%  "

$ws.Range("E21").Value = "  +0.53%  "

$ws.Range("E22").Value = "  +0.94%  "

$ws.Range("D23").Value = "9.18"
$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("E24").Value = "  -0.28%  "

$ws.Range("D25").Value = "154.02"
$ws.Range("E25").Value = "  +0.37%  "

$ws.Range("D26").Value = "6.59"
$ws.Range("E26").Value = "  -0.65%  "

$ws.Range("D27").Value = "15.03"
$ws.Range("E27").Value = "  +0.28%  "

$ws.Range("E28").Value = "  +1.13%  "

$ws.Range("E29").Value = "  +0.45%  "

$ws.Range("D30").Value = "1.13"
$ws.Range("E30").Value = "  +4.16%  "

$ws.Range("D31").Value = "0.0472"
$ws.Range("E31").Value = "  +0.39%  "

$ws.Range("E32").Value = "  +0.48%  "

$ws.Range("D33").Value = "3.18"
$ws.Range("E33").Value = "  +1.91%  "

$ws.Range("D34").Value = "1.429.69"
$ws.Range("E34").Value = "  +0.79%  "

$ws.Range("D35").Value = "1.09"
$ws.Range("E35").Value = "  +0.66%  "

$ws.Range("E36").Value = "  -0.65%  "

$ws.Range("E37").Value = "  +1.80%  "

$ws.Range("D38").Value = "0.0167"
$ws.Range("E38").Value = "  +0.86%  "

$ws.Range("D39").Value = "0.529"
$ws.Range("E39").Value = "  +0.45%  "

$ws.Range("D40").Value = "5.80"
$ws.Range("E40").Value = "  +2.73%  "

$ws.Range("D41").Value = "0.806"
$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("E42").Value = "  +0.52%  "

$ws.Range("D43").Value = "2.35"
$ws.Range("E43").Value = "  +1.35%  "

$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.28%  "

$ws.Range("D45").Value = "64.34"
$ws.Range("E45").Value = "  -0.23%  "

$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("D47").Value = "1.701.74"
$ws.Range("E47").Value = "  +0.90%  "

$ws.Range("D48").Value = "85.95"

$ws.Range("E49").Value = "  +2.56%  "

$ws.Range("D50").Value = "0.0518"
$ws.Range("E50").Value = "  -0.60%  "

$ws.Range("D51").Value = "0.0958"
$ws.Range("E51").Value = "  -0.14%  "
